# Apply crypto price/volume updates (GitHub Actions refresh simulation)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: assign directly
$ws.Range("D2").Value = '23.269.28'
$ws.Range("D3").Value = '1.599.00'
$ws.Range("E3").Value = '  -3.33%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("E6").Value = '  -2.31%  '
$ws.Range("E7").Value = '  -3.24%  '
$ws.Range("E8").Value = '  -4.61%  '
$ws.Range("E9").Value = '  -1.49%  '
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("E11").Value = '  -5.35%  '
$ws.Range("E12").Value = '  -3.93%  '
$ws.Range("E13").Value = '  -4.13%  '
$ws.Range("E14").Value = '  -7.05%  '
$ws.Range("E15").Value = '  -5.04%  '
$ws.Range("E16").Value = '  -2.52%  '
$ws.Range("D17").Value = '1.598.00'
$ws.Range("E17").Value = '  -3.55%  '
$ws.Range("E18").Value = '  -3.14%  '
$ws.Range("E19").Value = '  -2.42%  '
$ws.Range("E20").Value = '  -6.33%  '
$ws.Range("E21").Value = '  -3.94%  '
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("E23").Value = '  -3.29%  '
$ws.Range("D24").Value = '23.305.48'
$ws.Range("E24").Value = '  -2.36%  '
$ws.Range("E25").Value = '  -3.41%  '
$ws.Range("E26").Value = '  -2.40%  '
$ws.Range("E27").Value = '  -3.86%  '
$ws.Range("E28").Value = '  -1.36%  '
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("E30").Value = '  -4.82%  '
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("E32").Value = '  -5.56%  '
$ws.Range("D33").Value = '1.777.78'
$ws.Range("E33").Value = '  -3.29%  '
$ws.Range("E34").Value = '  -6.37%  '
$ws.Range("E35").Value = '  -4.24%  '
$ws.Range("E36").Value = '  -5.71%  '
$ws.Range("E37").Value = '  -4.51%  '
$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("E38").Value = '  -5.95%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E39").Value = '  -4.93%  '
$ws.Range("E40").Value = '  -2.70%  '
$ws.Range("E41").Value = '  -1.27%  '
$ws.Range("E42").Value = '  -3.93%  '
$ws.Range("E43").Value = '  -4.64%  '
$ws.Range("E44").Value = '  -1.11%  '
$ws.Range("E45").Value = '  -4.65%  '
$ws.Range("E46").Value = '  -5.44%  '
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("E48").Value = '  -2.15%  '
$ws.Range("E49").Value = '  -2.57%  '
$ws.Range("E50").Value = '  -2.01%  '
$ws.Range("E51").Value = '  -4.37%  '

# Numeric-looking text values: must stay as TEXT (not be converted to numbers),
# matching the source inline-string cells. Force text format, write, then restore
# the original style so no stray formatting is introduced.
$origStyle = $ws.Range("D2:D51").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '302.22'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3767'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3662'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '50.18'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.004'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.279'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08128'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.87'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.657'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.540'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001273'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.43'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06823'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.56'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.651'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.18'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.395'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.969'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.23'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.89'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.331'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.72'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.468'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.404'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9731'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07725'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02795'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.349'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.25'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2554'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08878'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.399'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7222'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.87'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.08'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6635'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.328'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.986'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08032'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.70'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.178'
$ws.Range("D2:D51").Style = $origStyle
